# Daily attendance processing - reorder "Recorded By" (column G) entries.
# Rule observed from the target diff: within each comma-separated list of
# recorder names, entries that are NOT an email address (e.g. "System",
# "system") keep their original relative order and are moved to the front;
# entries that ARE an email address are sorted alphabetically
# (case-insensitive) and appended after them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -notmatch ",") { continue }

    $parts = $val -split "," | ForEach-Object { $_.Trim() }

    $nonEmail = @()
    $email = @()
    foreach ($p in $parts) {
        if ($p -like "*@*") {
            $email += $p
        } else {
            $nonEmail += $p
        }
    }

    $emailSorted = $email | Sort-Object { $_.ToLower() }

    $newParts = @()
    $newParts += $nonEmail
    $newParts += $emailSorted

    $newVal = $newParts -join ", "

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
